$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 319266.66
$ws.Range("I28").Value = 411934.34
$ws.Range("J28").Value = 6513.25
$ws.Range("K28").Value = 411934.34
$ws.Range("L28").Value = 6513.25
$ws.Range("M28").Value = -411449.34
$ws.Range("N28").Value = -7483.25

$ws.Range("H98").Value = 560989.4
$ws.Range("I98").Value = 623099.3
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 623099.3
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -621601.3
$ws.Range("N98").Value = -4996

$ws.Range("H116").Value = 3877.7778
$ws.Range("I116").Value = 4271.4287
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 4271.4287
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -829.4287000000004
$ws.Range("N116").Value = -9384

$ws.Range("H122").Value = 560989.4
$ws.Range("I122").Value = 623099.3
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 1869297.9
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1866847.9
$ws.Range("N122").Value = -10900

$ws.Range("H138").Value = 9859657
$ws.Range("I138").Value = 7148094.5
$ws.Range("J138").Value = 10424566
$ws.Range("K138").Value = 21444283.5
$ws.Range("L138").Value = 31273698
$ws.Range("M138").Value = -21439143.5
$ws.Range("N138").Value = -31283978

$ws.Range("H141").Value = 2714.3125
$ws.Range("I141").Value = 2685.75
$ws.Range("J141").Value = 2800
$ws.Range("K141").Value = 8057.25
$ws.Range("L141").Value = 8400
$ws.Range("M141").Value = -2877.25
$ws.Range("N141").Value = -18760

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 92778.63
$ws.Range("I2").Value = 101906.5
$ws.Range("K2").Value = 101906.5
$ws.Range("M2").Value = -101793.5

$ws.Range("H32").Value = 70345.52
$ws.Range("I32").Value = 15326
$ws.Range("J32").Value = 105715.21
$ws.Range("K32").Value = 15326
$ws.Range("L32").Value = 105715.21
$ws.Range("M32").Value = -15039
$ws.Range("N32").Value = -106289.21

$ws.Range("H45").Value = 1491.8636
$ws.Range("I45").Value = 1442.7646
$ws.Range("J45").Value = 1658.8
$ws.Range("K45").Value = 1442.7646
$ws.Range("L45").Value = 1658.8
$ws.Range("M45").Value = -1065.7646
$ws.Range("N45").Value = -2412.8

$ws.Range("H61").Value = 2686.3
$ws.Range("I61").Value = 2285.9412
$ws.Range("J61").Value = 4955
$ws.Range("K61").Value = 2285.9412
$ws.Range("L61").Value = 4955
$ws.Range("M61").Value = -2073.9412
$ws.Range("N61").Value = -5379

$ws.Range("H63").Value = 22599.375
$ws.Range("I63").Value = 24427.857
$ws.Range("K63").Value = 24427.857
$ws.Range("M63").Value = -23741.857

$ws.Range("H66").Value = 22599.375
$ws.Range("I66").Value = 24427.857
$ws.Range("K66").Value = 122139.285
$ws.Range("M66").Value = -118707.285

$ws.Range("H116").Value = 92778.63
$ws.Range("I116").Value = 101906.5
$ws.Range("K116").Value = 101906.5
$ws.Range("M116").Value = -99612.5

$ws.Range("H122").Value = 1202.9445
$ws.Range("I122").Value = 1243.4706
$ws.Range("K122").Value = 3730.4118
$ws.Range("M122").Value = -1280.4118

$ws.Range("H136").Value = 2686.3
$ws.Range("I136").Value = 2285.9412
$ws.Range("J136").Value = 4955
$ws.Range("K136").Value = 6857.823600000001
$ws.Range("L136").Value = 14865
$ws.Range("M136").Value = -4307.823600000001
$ws.Range("N136").Value = -19965

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 92778.63
$ws.Range("I3").Value = 101906.5
$ws.Range("K3").Value = 101906.5
$ws.Range("M3").Value = -101792.5

$ws.Range("H132").Value = 45668.332
$ws.Range("J132").Value = 45668.332
$ws.Range("L132").Value = 45668.332
$ws.Range("N132").Value = -55788.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5419
$ws.Range("I31").Value = 2750.7334
$ws.Range("J31").Value = 7920.5
$ws.Range("K31").Value = 2750.7334
$ws.Range("L31").Value = 7920.5
$ws.Range("M31").Value = -2455.7334
$ws.Range("N31").Value = -8510.5

$ws.Range("H34").Value = 5419
$ws.Range("I34").Value = 2750.7334
$ws.Range("J34").Value = 7920.5
$ws.Range("K34").Value = 2750.7334
$ws.Range("L34").Value = 7920.5
$ws.Range("M34").Value = -2548.7334
$ws.Range("N34").Value = -8324.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3813
$ws.Range("I51").Value = 704
$ws.Range("J51").Value = 4257.143
$ws.Range("K51").Value = 2112
$ws.Range("L51").Value = 12771.429
$ws.Range("M51").Value = -1652
$ws.Range("N51").Value = -13691.429

$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 7000
$ws.Range("L55").Value = 21000
$ws.Range("N55").Value = -21354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1793.9333
$ws.Range("I97").Value = 1663.75
$ws.Range("J97").Value = 1942.7142
$ws.Range("K97").Value = 1663.75
$ws.Range("L97").Value = 1942.7142
$ws.Range("M97").Value = -1167.75
$ws.Range("N97").Value = -2934.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3545.3635
$ws.Range("J7").Value = 3500.2222
$ws.Range("L7").Value = 3500.2222
$ws.Range("N7").Value = -3724.2222

$ws.Range("H22").Value = 13780.25
$ws.Range("I22").Value = 1690
$ws.Range("J22").Value = 33930.668
$ws.Range("K22").Value = 1690
$ws.Range("L22").Value = 33930.668
$ws.Range("M22").Value = -1395
$ws.Range("N22").Value = -34520.668

$ws.Range("H27").Value = 13780.25
$ws.Range("I27").Value = 1690
$ws.Range("J27").Value = 33930.668
$ws.Range("K27").Value = 1690
$ws.Range("L27").Value = 33930.668
$ws.Range("M27").Value = -1583
$ws.Range("N27").Value = -34144.668

$ws.Range("H55").Value = 866.6667
$ws.Range("I55").Value = 856.8570999999999
$ws.Range("J55").Value = 901
$ws.Range("K55").Value = 856.8570999999999
$ws.Range("L55").Value = 901
$ws.Range("M55").Value = -683.8570999999999
$ws.Range("N55").Value = -1247

$ws.Range("H100").Value = 1118208.6
$ws.Range("I100").Value = 2605828.5
$ws.Range("J100").Value = 2493.75
$ws.Range("K100").Value = 2605828.5
$ws.Range("L100").Value = 2493.75
$ws.Range("M100").Value = -2605287.5
$ws.Range("N100").Value = -3575.75

$ws.Range("H126").Value = 3545.3635
$ws.Range("J126").Value = 3500.2222
$ws.Range("L126").Value = 10500.6666
$ws.Range("N126").Value = -15440.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1487.0834
$ws.Range("I81").Value = 834.3
$ws.Range("J81").Value = 4751
$ws.Range("K81").Value = 1668.6
$ws.Range("L81").Value = 9502
$ws.Range("M81").Value = -607.5999999999999
$ws.Range("N81").Value = -11624

$ws.Range("H84").Value = 1487.0834
$ws.Range("I84").Value = 834.3
$ws.Range("J84").Value = 4751
$ws.Range("K84").Value = 8343
$ws.Range("L84").Value = 47510
$ws.Range("M84").Value = -3039
$ws.Range("N84").Value = -58118

$ws.Range("H122").Value = 995
$ws.Range("I122").Value = 995
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2985
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -535
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 28207.639
$ws.Range("I126").Value = 37383.965
$ws.Range("K126").Value = 112151.895
$ws.Range("M126").Value = -109681.895
